# Update Name of Algo
# Apply updated numeric results to the KNN imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C11"  = -12.274
    "A12"  = -21.451
    "C23"  = -12.309
    "D24"  = -7.601999999999999
    "C28"  = -12.881
    "A32"  = -21.977
    "C32"  = -13.538
    "C34"  = -12.101
    "A36"  = -20.43
    "A38"  = -20.03
    "D38"  = -8.279
    "C42"  = -12.36
    "A46"  = -21.744
    "D52"  = -7.679
    "A54"  = -21.832
    "C54"  = -13.054
    "A55"  = -22.016
    "A67"  = -21.577
    "A69"  = -21.47
    "A72"  = -21.689
    "D78"  = -7.961999999999999
    "D83"  = -7.898999999999999
    "D85"  = -8.664000000000001
    "D86"  = -8.312000000000001
    "A91"  = -20.755
    "D96"  = -7.444
    "C97"  = -11.509
    "A99"  = -21.785
    "C99"  = -12.302
    "C101" = -12.188
    "D103" = -8.315999999999999
    "A104" = -21.175
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
